$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted ahead of the existing "Papa" rows
# for Vega Monumental Concepción (Bíobío), pushing the former rows 364-379
# down to 365-380 while keeping their data intact.
$ws.Rows.Item(364).Insert()

# Fill in the newly inserted row 364 with the new record's data.
$ws.Range("A364").Value = 11
$ws.Range("B364").Value = "Vega Monumental Concepción"
$ws.Range("C364").Value = "Bíobío"
$ws.Range("D364").Value2 = 44939
$ws.Range("E364").Value = 8
$ws.Range("F364").Value = 100114001
$ws.Range("G364").Value = "Papa"
$ws.Range("H364").Value = "Asterix"
$ws.Range("I364").Value = "1a (cosecha)"
$ws.Range("J364").Value = 5000
$ws.Range("K364").Value = 12000
$ws.Range("L364").Value = 12500
$ws.Range("M364").Value = 12250
$ws.Range("N364").Value = "$/saco 25 kilos"
$ws.Range("O364").Value = "Región de La Araucanía"
$ws.Range("P364").Value = 490
$ws.Range("Q364").Value = 25
$ws.Range("R364").Value = "Hortaliza"
